$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D2:D51) keeps its original text formatting so that
# numeric-looking values (e.g. "0.999", "7.03") are stored as text, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '64.496.19'
$ws.Range('E2').Value = '  +4.96%  '
$ws.Range('D3').Value = '3.099.37'
$ws.Range('E3').Value = '  +3.80%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '558.44'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '143.96'
$ws.Range('E6').Value = '  +10.20%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.099.34'
$ws.Range('E8').Value = '  +3.93%  '
$ws.Range('D9').Value = '0.500'
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').Value = '7.03'
$ws.Range('E10').Value = '  +17.35%  '
$ws.Range('E11').Value = '  +5.80%  '
$ws.Range('D12').Value = '0.461'
$ws.Range('E12').Value = '  +4.09%  '
$ws.Range('D13').Value = '0.0000227'
$ws.Range('E13').Value = '  +4.28%  '
$ws.Range('D14').Value = '35.08'
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('D15').Value = '3.598.06'
$ws.Range('E15').Value = '  +3.93%  '
$ws.Range('D16').Value = '64.550.88'
$ws.Range('E16').Value = '  +4.89%  '
$ws.Range('D17').Value = '3.093.75'
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').Value = '6.72'
$ws.Range('D20').Value = '484.29'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').Value = '13.75'
$ws.Range('E21').Value = '  +5.30%  '
$ws.Range('D22').Value = '0.671'
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').Value = '7.55'
$ws.Range('E23').Value = '  +8.46%  '
$ws.Range('D24').Value = '13.22'
$ws.Range('E24').Value = '  +10.45%  '
$ws.Range('D25').Value = '80.94'
$ws.Range('E25').Value = '  +0.87%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '2.82'
$ws.Range('E27').Value = '  +4.18%  '
$ws.Range('D28').Value = '8.00'
$ws.Range('E28').Value = '  +4.76%  '
$ws.Range('D29').Value = '2.06'
$ws.Range('E29').Value = '  +8.82%  '
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').Value = '26.09'
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('E32').Value = '  +2.51%  '
$ws.Range('D33').Value = '2.44'
$ws.Range('E33').Value = '  +6.06%  '
$ws.Range('D34').Value = '5.76'
$ws.Range('E34').Value = '  +5.30%  '
$ws.Range('E35').Value = '  +1.35%  '
$ws.Range('D36').Value = '6.12'
$ws.Range('E36').Value = '  +4.94%  '
$ws.Range('D37').Value = '463.42'
$ws.Range('E37').Value = '  +3.32%  '
$ws.Range('D38').Value = '0.0407'
$ws.Range('E38').Value = '  +7.34%  '
$ws.Range('D39').Value = '0.0826'
$ws.Range('E39').Value = '  +4.90%  '
$ws.Range('D40').Value = '3.025.74'
$ws.Range('E40').Value = '  -2.70%  '
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('D42').Value = '8.30'
$ws.Range('E42').Value = '  +3.14%  '
$ws.Range('D43').Value = '2.70'
$ws.Range('E43').Value = '  +16.33%  '
$ws.Range('D44').Value = '28.04'
$ws.Range('E44').Value = '  +10.72%  '
$ws.Range('D45').Value = '0.260'
$ws.Range('E45').Value = '  +8.09%  '
$ws.Range('D47').Value = '2.09'
$ws.Range('E47').Value = '  +8.02%  '
$ws.Range('E48').Value = '  +4.64%  '
$ws.Range('D49').Value = '118.60'
$ws.Range('E49').Value = '  +4.26%  '
$ws.Range('D50').Value = '0.0₃0517'
$ws.Range('E50').Value = '  +8.12%  '
$ws.Range('D51').Value = '2.06'
$ws.Range('E51').Value = '  +3.34%  '
